$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 2185.2856
$ws.Range("I38").Value = 534.1667
$ws.Range("J38").Value = 3423.625
$ws.Range("K38").Value = 1602.5001
$ws.Range("L38").Value = 10270.875
$ws.Range("M38").Value = -1230.5001
$ws.Range("N38").Value = -11014.875
$ws.Range("H43").Value = 2999.6667
$ws.Range("I43").Value = 2999.5
$ws.Range("K43").Value = 2999.5
$ws.Range("M43").Value = -2930.5
$ws.Range("H104").Value = 982
$ws.Range("I104").Value = 982
$ws.Range("K104").Value = 2946
$ws.Range("M104").Value = -1199
$ws.Range("H129").Value = 1950.5385
$ws.Range("J129").Value = 2407.375
$ws.Range("L129").Value = 7222.125
$ws.Range("N129").Value = -17222.125
$ws.Range("H138").Value = 4429.724
$ws.Range("I138").Value = 2110.7778
$ws.Range("K138").Value = 6332.3334
$ws.Range("M138").Value = -1192.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 3284.5
$ws.Range("I19").Value = 3802
$ws.Range("J19").Value = 2249.5
$ws.Range("K19").Value = 3802
$ws.Range("L19").Value = 2249.5
$ws.Range("M19").Value = -3573
$ws.Range("N19").Value = -2707.5
$ws.Range("H63").Value = 4833.3335
$ws.Range("I63").Value = 4000
$ws.Range("J63").Value = 5250
$ws.Range("K63").Value = 4000
$ws.Range("L63").Value = 5250
$ws.Range("M63").Value = -3314
$ws.Range("N63").Value = -6622
$ws.Range("H66").Value = 4833.3335
$ws.Range("I66").Value = 4000
$ws.Range("J66").Value = 5250
$ws.Range("K66").Value = 20000
$ws.Range("L66").Value = 26250
$ws.Range("M66").Value = -16568
$ws.Range("N66").Value = -33114
$ws.Range("H132").Value = 737.8
$ws.Range("I132").Value = 737.8
$ws.Range("K132").Value = 2213.4
$ws.Range("M132").Value = 316.6000000000004
$ws.Range("H135").Value = 86952.336
$ws.Range("J135").Value = 86952.336
$ws.Range("L135").Value = 86952.336
$ws.Range("N135").Value = -97092.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 146.66667
$ws.Range("I11").Value = 400
$ws.Range("J11").Value = 20
$ws.Range("K11").Value = 400
$ws.Range("L11").Value = 20
$ws.Range("M11").Value = -260
$ws.Range("N11").Value = -300
$ws.Range("H14").Value = 812
$ws.Range("I14").Value = 650
$ws.Range("J14").Value = 866
$ws.Range("K14").Value = 650
$ws.Range("L14").Value = 866
$ws.Range("M14").Value = -478
$ws.Range("N14").Value = -1210
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("N23").ClearContents()
$ws.Range("H31").Value = 5025
$ws.Range("J31").Value = 5025
$ws.Range("L31").Value = 5025
$ws.Range("N31").Value = -5529

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 689.3
$ws.Range("I2").Value = 550
$ws.Range("J2").Value = 782.1667
$ws.Range("K2").Value = 550
$ws.Range("L2").Value = 782.1667
$ws.Range("M2").Value = -437
$ws.Range("N2").Value = -1008.1667
$ws.Range("H5").Value = 733.8
$ws.Range("I5").Value = 262.7143
$ws.Range("J5").Value = 1833
$ws.Range("K5").Value = 262.7143
$ws.Range("L5").Value = 1833
$ws.Range("M5").Value = -150.7143
$ws.Range("N5").Value = -2057
$ws.Range("H12").Value = 531.1429000000001
$ws.Range("I12").Value = 767
$ws.Range("K12").Value = 767
$ws.Range("M12").Value = -597
$ws.Range("H14").Value = 4309.6
$ws.Range("J14").Value = 4500
$ws.Range("L14").Value = 4500
$ws.Range("N14").Value = -4840
$ws.Range("H31").Value = 5614.727
$ws.Range("I31").Value = 2343
$ws.Range("K31").Value = 2343
$ws.Range("M31").Value = -2048
$ws.Range("H34").Value = 5614.727
$ws.Range("I34").Value = 2343
$ws.Range("K34").Value = 2343
$ws.Range("M34").Value = -2141
$ws.Range("H52").Value = 191295.72
$ws.Range("J52").Value = 215680
$ws.Range("L52").Value = 215680
$ws.Range("N52").Value = -216268

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 24.652174
$ws.Range("I2").Value = 27
$ws.Range("J2").Value = 21
$ws.Range("K2").Value = 162
$ws.Range("L2").Value = 126
$ws.Range("M2").Value = -49
$ws.Range("N2").Value = -352
$ws.Range("H9").Value = 184.25
$ws.Range("I9").Value = 79
$ws.Range("J9").Value = 500
$ws.Range("K9").Value = 237
$ws.Range("L9").Value = 1500
$ws.Range("M9").Value = -13
$ws.Range("N9").Value = -1948
$ws.Range("H44").Value = 821.75
$ws.Range("I44").Value = 194.57143
$ws.Range("K44").Value = 583.71429
$ws.Range("M44").Value = -185.71429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 153
$ws.Range("I9").Value = 186.85715
$ws.Range("J9").Value = 93.75
$ws.Range("K9").Value = 186.85715
$ws.Range("L9").Value = 93.75
$ws.Range("M9").Value = -16.85714999999999
$ws.Range("N9").Value = -433.75
$ws.Range("H13").Value = 232.25
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 232.25
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 232.25
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -510.25
$ws.Range("H17").Value = 50000
$ws.Range("J17").Value = 50000
$ws.Range("L17").Value = 50000
$ws.Range("N17").Value = -50336
$ws.Range("H25").Value = 1706.6666
$ws.Range("J25").Value = 560
$ws.Range("L25").Value = 560
$ws.Range("N25").Value = -1618
$ws.Range("H132").Value = 64483.25
$ws.Range("I132").Value = 73102.36
$ws.Range("J132").Value = 4149.5
$ws.Range("K132").Value = 219307.08
$ws.Range("L132").Value = 12448.5
$ws.Range("M132").Value = -216777.08
$ws.Range("N132").Value = -17508.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H132").Value = 10363
$ws.Range("I132").Value = 3223.5
$ws.Range("K132").Value = 9670.5
$ws.Range("M132").Value = -7140.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 8818682
$ws.Range("J5").Value = 4500083.5
$ws.Range("L5").Value = 4500083.5
$ws.Range("N5").Value = -4500307.5
$ws.Range("H14").Value = 234.48387
$ws.Range("I14").Value = 312.77777
$ws.Range("K14").Value = 312.77777
$ws.Range("M14").Value = -144.77777
$ws.Range("H100").Value = 1242.1818
$ws.Range("I100").Value = 1473.7778
$ws.Range("K100").Value = 2947.5556
$ws.Range("M100").Value = -2406.5556
$ws.Range("H126").Value = 4897.727
$ws.Range("I126").Value = 1118.75
$ws.Range("J126").Value = 7057.143
$ws.Range("K126").Value = 3356.25
$ws.Range("L126").Value = 21171.429
$ws.Range("M126").Value = -886.25
$ws.Range("N126").Value = -26111.429
$ws.Range("H132").Value = 1406.6364
$ws.Range("I132").Value = 1406.6364
$ws.Range("K132").Value = 4219.9092
$ws.Range("M132").Value = -1689.9092
